$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit cyclically rotates the per-record data (Id, Taxonsorteringsordning,
# Rödlistade, TaxonId, Artnamn, Vetenskapligt namn, Auktor, Kön, Ost, Nord,
# Noggrannhet) among rows 10-14, and swaps the same data between rows 18-19.
# All other columns (location, county, dates, observers, ...) are identical
# across these rows already, so only the columns below need updating.

# --- Row 10 (becomes former row 14's "Lunglav" record) ---
$ws.Range("A10").Value = 112129990
$ws.Range("B10").Value = 78713
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 6458
$ws.Range("F10").Value = "Lunglav"
$ws.Range("G10").Value = "Lobaria pulmonaria"
$ws.Range("H10").Value = "(L.) Hoffm."
$ws.Range("L10").ClearContents()
$ws.Range("Q10").Value = 489779
$ws.Range("R10").Value = 6951289
$ws.Range("S10").Value = 1

# --- Row 11 (becomes former row 10's "Knärot" record) ---
$ws.Range("A11").Value = 112129159
$ws.Range("B11").Value = 96735
$ws.Range("D11").Value = "VU"
$ws.Range("E11").Value = 220787
$ws.Range("F11").Value = "Knärot"
$ws.Range("G11").Value = "Goodyera repens"
$ws.Range("H11").Value = "(L.) R. Br."
$ws.Range("Q11").Value = 489738
$ws.Range("R11").Value = 6951149
$ws.Range("S11").Value = 5

# --- Row 12 (becomes former row 11's "Knärot" record) ---
$ws.Range("A12").Value = 112128360
$ws.Range("B12").Value = 96735
$ws.Range("D12").Value = "VU"
$ws.Range("E12").Value = 220787
$ws.Range("F12").Value = "Knärot"
$ws.Range("G12").Value = "Goodyera repens"
$ws.Range("H12").Value = "(L.) R. Br."
$ws.Range("Q12").Value = 489808
$ws.Range("R12").Value = 6951101
$ws.Range("S12").Value = 5

# --- Row 13 (becomes former row 12's "Knärot" record) ---
$ws.Range("A13").Value = 112128076
$ws.Range("B13").Value = 96735
$ws.Range("D13").Value = "VU"
$ws.Range("E13").Value = 220787
$ws.Range("F13").Value = "Knärot"
$ws.Range("G13").Value = "Goodyera repens"
$ws.Range("H13").Value = "(L.) R. Br."
$ws.Range("Q13").Value = 489837
$ws.Range("R13").Value = 6951074
$ws.Range("S13").Value = 1

# --- Row 14 (becomes former row 13's "Knärot" record) ---
$ws.Range("A14").Value = 112130763
$ws.Range("B14").Value = 96735
$ws.Range("D14").Value = "VU"
$ws.Range("E14").Value = 220787
$ws.Range("F14").Value = "Knärot"
$ws.Range("G14").Value = "Goodyera repens"
$ws.Range("H14").Value = "(L.) R. Br."
$ws.Range("L14").ClearContents()
$ws.Range("L14").Style = "Normal"
$ws.Range("Q14").Value = 489727
$ws.Range("R14").Value = 6951335
$ws.Range("S14").Value = 5

# --- Rows 18 & 19 swap their Id / Ost / Nord / Noggrannhet values ---
$ws.Range("A18").Value = 112130099
$ws.Range("Q18").Value = 489754
$ws.Range("R18").Value = 6951300
$ws.Range("S18").Value = 1

$ws.Range("A19").Value = 112128285
$ws.Range("Q19").Value = 489826
$ws.Range("R19").Value = 6951101
$ws.Range("S19").Value = 5
